# Widgets - Date picker
# Adds locator rows for the "Auto Complete" and "Date Picker" widgets
# to the locators repository sheet, and updates the sheet view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 135: "Auto Complete" section header -----------------------------
$ws.Range("A135").Value = "Auto Complete"

# --- Row 136: nav locator --------------------------------------------------
$ws.Range("A136").Value = "autoCompleteNav"
$ws.Range("B136").Value = "//div[@class='element-list collapse show']//li[@id='item-1']"
$ws.Range("C136").Value = "By.xpath"

# --- Row 137: input locator -------------------------------------------------
$ws.Range("A137").Value = "autoCompleteInput"
$ws.Range("B137").Value = '//*[@id="autoCompleteMultipleContainer"]'
$ws.Range("C137").Value = "By.xpath"

# --- Row 138: suggestion locator --------------------------------------------
$ws.Range("A138").Value = "suggestion1"
$ws.Range("B138").Value = '//*[@id="react-select-8-option-0"]'
$ws.Range("C138").Value = "By.xpath"

# --- Row 139: input (selected value) locator --------------------------------
$ws.Range("A139").Value = "autoCompleteInput"
$ws.Range("B139").Value = '//*[@id="autoCompleteMultipleContainer"]/div/div[1]/div[1]/div[1]'
$ws.Range("C139").Value = "By.xpath"

# --- Row 140: "Date Picker" section header -----------------------------------
$ws.Range("A140").Value = "datePicker"

# --- Row 141: nav locator -----------------------------------------------------
$ws.Range("A141").Value = "datePickerNav"
$ws.Range("B141").Value = "//div[@class='element-list collapse show']//li[@id='item-2']"
$ws.Range("C141").Value = "By.xpath"

# --- Row 142: widgets scroll locator -------------------------------------------
$ws.Range("A142").Value = "WidgetsScroll"
$ws.Range("B142").Value = "//body/div[@id='app']/div[@class='body-height']/div[@class='container playgound-body']/div[@class='row']/div[1]/div[1]/div[1]/div[4]/span[1]/div[1]"
$ws.Range("C142").Value = "By.xpath"

# --- Row 143: date picker input locator -----------------------------------------
$ws.Range("A143").Value = "datePickerInput"
$ws.Range("B143").Value = '//*[@id="datePickerMonthYear"]/div'
$ws.Range("C143").Value = "By.xpath"

# --- Row 144: date picker selected date locator ----------------------------------
$ws.Range("A144").Value = "datePickerSelectedDate"
$ws.Range("B144").Value = '//*[@id="datePickerMonthYearInput"]'
$ws.Range("C144").Value = "By.xpath"

# --- Update the view: scroll position and active selection -----------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 136
$win.ScrollColumn = 1
$ws.Range("B148").Select()
